$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "count" values for existing rows 2-34 (column B)
$updatedCounts = @(
    1771, 1792, 1830, 1840, 1847, 1843, 1837, 1838, 1839, 1850,
    1815, 1784, 1786, 1817, 1893, 1739, 1676, 1707, 1742, 1726,
    1727, 1739, 1751, 1736, 1643, 1657, 1714, 1760, 1824, 1823,
    1822, 1816, 1838
)

for ($i = 0; $i -lt $updatedCounts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $updatedCounts[$i]
}

# New rows 35-42: date (ts), count (B), time_unit (C) = "M"
$newRows = @(
    @(44500, 1824),
    @(44530, 1821),
    @(44561, 1783),
    @(44592, 1730),
    @(44620, 1832),
    @(44651, 1875),
    @(44681, 1913),
    @(44712, 1940)
)

$startRow = 35
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $startRow + $i
    $dateSerial = $newRows[$i][0]
    $countVal = $newRows[$i][1]

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $dateSerial
    $aCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $countVal
    $ws.Cells.Item($row, 3).Value = "M"
}
